$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "This program will be a very basic chat bot. Essentially it will..." --
#    remove the gramStart/gramEnd proofErr markers around "Essentially" and
#    collapse the three runs they split into one plain run.
# ---------------------------------------------------------------------------
$emDash = [char]0x2014
$para1 = "bot. Essentially it will simulate very basic Artificial Intelligence (AI) through interaction with the user. The program will accept input in the form of grammatically correct English and respond appropriately" + $emDash + "also in English."

$rng = $d.Content
$rng.Find.Execute($para1, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$start1 = $rng.Start
$end1 = $rng.End
$full1 = $d.Range($start1, $end1)
$full1.Delete()
$ins1 = $d.Range($start1, $start1)
$ins1.InsertAfter($para1)

# ---------------------------------------------------------------------------
# 2) "The program will utilize dictionary files ... (i.e. noun, verb, ...)" --
#    remove the gramStart/gramEnd proofErr markers around "i.e" and collapse
#    the three runs they split into one plain run.
# ---------------------------------------------------------------------------
$para2 = "The program will utilize dictionary files to detect words and determine the part of speech they belong to. (i.e. noun, verb, article, etc.) It will then attempt to construct a comprehensible sentence by analyzing how those parts of speech should interact based on a set of rules."

$rng2 = $d.Content
$rng2.Find.Execute($para2, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$start2 = $rng2.Start
$end2 = $rng2.End
$full2 = $d.Range($start2, $end2)
$full2.Delete()
$ins2 = $d.Range($start2, $start2)
$ins2.InsertAfter($para2)

# ---------------------------------------------------------------------------
# 3) Limitations paragraph -- type "truly " in front of "conjugating words"
#    and let the _GoBack bookmark (Word's "last edit" marker) follow the new
#    text, moving it away from the end of the document.
# ---------------------------------------------------------------------------
$rng3 = $d.Content
$rng3.Find.Execute("conjugating words", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$posConj = $rng3.Start

$insTruly = $d.Range($posConj, $posConj)
$insTruly.InsertAfter("truly ")

# Re-locate "conjugating words" now that "truly " sits in front of it.
$rng4 = $d.Content
$rng4.Find.Execute("conjugating words", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$posConjNew = $rng4.Start

# Re-seat the _GoBack bookmark right before "conjugating" -- adding a
# bookmark with that reserved name automatically removes it from wherever
# it used to live (end of the document, after "Steve Halladay").
$goBackRange = $d.Range($posConjNew, $posConjNew)
$d.Bookmarks.Add("_GoBack", $goBackRange)

# Split "truly " off into its own run (matching how Word leaves a freshly
# typed word in its own run) by briefly bookmarking the point just before it
# and then discarding that scratch bookmark -- the run boundary it created
# survives the bookmark's removal.
$posBeforeTruly = $posConjNew - "truly ".Length
$scratchRange = $d.Range($posBeforeTruly, $posBeforeTruly)
$d.Bookmarks.Add("zzzScratchSplit", $scratchRange)
$d.Bookmarks.Item("zzzScratchSplit").Delete()

Write-Output "Edit complete"
